$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Telavi")

# Extend the left table one more year (2023) by copying the formatting
# from the last existing year column (J) into the new column (K),
# then filling in the 2023 figures.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1385.7
$ws.Range("K5").Value = 838.8
$ws.Range("K6").Value = 1897.9
